# Generate Report for Handoff
#
# For every row whose Status is "Ready for handoff" in the zh-cn and de-de
# sheets, mark the Priority as handed-off ("ht") and stamp a fresh
# "Latest Handoff Datetime".

$wb = $excel.ActiveWorkbook

$sheets = @(
    @{ Name = "zh-cn"; Stamp = "2016-08-13 08:34:35" },
    @{ Name = "de-de"; Stamp = "2016-08-13 08:34:42" }
)

foreach ($sheetInfo in $sheets) {
    $ws = $wb.Worksheets.Item($sheetInfo.Name)

    for ($row = 2; $row -le 7; $row++) {
        $status = $ws.Cells.Item($row, 3).Value2
        if ($status -eq "Ready for handoff") {
            $ws.Cells.Item($row, 5).Value = "ht"
            $ws.Cells.Item($row, 8).Value = $sheetInfo.Stamp
        }
    }
}
